{"js": "// Remove the `{{ p.qr_url }}` merge-field placeholder text from the QR-code\n// table cell, leaving the paragraph (and its table cell) empty.\nconst body = context.document.body;\n\n// Locate the merge-field text wherever it lives in the document body\n// (it sits in the last column of the second table, a vMerge-restart cell).\nconst results = body.search(\"{{ p.qr_url }}\", { matchCase: true, matchWildcards: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the '{{ p.qr_url }}' placeholder text.\");\n}\n\n// Clear the whole paragraph's content (runs + proofErr markers) while\n// keeping the empty <w:p/> paragraph mark in place, matching how Word's\n// Range.Clear / paragraph.clear() behave.\nconst hitParagraph = results.items[0].paragraphs.getFirst();\nhitParagraph.clear();\n\nawait context.sync();\n", "ps1": "# Remove the `{{ p.qr_url }}` merge-field placeholder text from the QR-code\n# table cell, leaving an empty paragraph behind (the <w:p> stays, its runs\n# and spell-check proofErr markers go away).\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n\n# MatchCase:=True, Forward:=True, Wrap:=wdFindContinue(1), Replace:=wdReplaceAll(2),\n# ReplaceWith:=\"\" -- a true Find & Replace removes every run that makes up\n# the matched text (unlike a plain `Range.Text = \"\"`, which only touches the\n# first run), so the paragraph mark survives empty.\n$found = $rng.Find.Execute(\"{{ p.qr_url }}\", $true, $false, $false, $false, $false, $true, 1, $false, \"\", 2)\n\nif (-not $found) {\n    throw \"Could not find the '{{ p.qr_url }}' placeholder text.\"\n}\n"}
